# "Generate Report for Archive"
# The localization status of the two handed-off files moved from
# "Ready for handoff" to "In Translation". Update every cell that carries
# that status (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3) and then re-fit the
# columns that were sized to the old (longer) status text - the shorter
# string lets those columns shrink.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = "In Translation"
$overview.Columns.Item(5).EntireColumn.AutoFit()
$overview.Columns.Item(6).EntireColumn.AutoFit()
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = "In Translation"
$zhcn.Columns.Item(3).EntireColumn.AutoFit()
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = "In Translation"
$dede.Columns.Item(3).EntireColumn.AutoFit()
$dede.Columns.Item(3).ColumnWidth = 12.5
